# Lattice-multiplication exercise table: refresh all 15 problem cells
# (5 rows x 3 columns) with the new set of factors / lattice digits,
# matching the regenerated output at commit 503736d.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11  # vertical-tab == <w:br/> line break inside a Range.Text assignment

$t.Cell(1,1).Range.Text = "15 x 67" + $nl + "  6    7" + $nl + "  ----" + $nl + "1|    |" + $nl + "5|    |"
$t.Cell(1,2).Range.Text = "67 x 30" + $nl + "  3    0" + $nl + "  ----" + $nl + "6|    |" + $nl + "7|    |"
$t.Cell(1,3).Range.Text = "17 x 19" + $nl + "  1    9" + $nl + "  ----" + $nl + "1|    |" + $nl + "7|    |"
$t.Cell(2,1).Range.Text = "90 x 93" + $nl + "  9    3" + $nl + "  ----" + $nl + "9|    |" + $nl + "0|    |"
$t.Cell(2,2).Range.Text = "39 x 37" + $nl + "  3    7" + $nl + "  ----" + $nl + "3|    |" + $nl + "9|    |"
$t.Cell(2,3).Range.Text = "65 x 74" + $nl + "  7    4" + $nl + "  ----" + $nl + "6|    |" + $nl + "5|    |"
$t.Cell(3,1).Range.Text = "15 x 27" + $nl + "  2    7" + $nl + "  ----" + $nl + "1|    |" + $nl + "5|    |"
$t.Cell(3,2).Range.Text = "76 x 23" + $nl + "  2    3" + $nl + "  ----" + $nl + "7|    |" + $nl + "6|    |"
$t.Cell(3,3).Range.Text = "72 x 87" + $nl + "  8    7" + $nl + "  ----" + $nl + "7|    |" + $nl + "2|    |"
$t.Cell(4,1).Range.Text = "24 x 82" + $nl + "  8    2" + $nl + "  ----" + $nl + "2|    |" + $nl + "4|    |"
$t.Cell(4,2).Range.Text = "69 x 58" + $nl + "  5    8" + $nl + "  ----" + $nl + "6|    |" + $nl + "9|    |"
$t.Cell(4,3).Range.Text = "61 x 41" + $nl + "  4    1" + $nl + "  ----" + $nl + "6|    |" + $nl + "1|    |"
$t.Cell(5,1).Range.Text = "61 x 51" + $nl + "  5    1" + $nl + "  ----" + $nl + "6|    |" + $nl + "1|    |"
$t.Cell(5,2).Range.Text = "57 x 97" + $nl + "  9    7" + $nl + "  ----" + $nl + "5|    |" + $nl + "7|    |"
$t.Cell(5,3).Range.Text = "16 x 74" + $nl + "  7    4" + $nl + "  ----" + $nl + "1|    |" + $nl + "6|    |"

Write-Host "Updated" $t.Rows.Count "rows x" $t.Columns.Count "columns of lattice problems."
